# Auto-generated edit script applying numeric corrections to Kraken_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 5 (G5=5503)
$ws.Range("H5").Value = 150.09091
$ws.Range("I5").Value = 155.1
$ws.Range("K5").Value = 155.1
$ws.Range("M5").Value = -40.09999999999999

# Row 12 (G12=5515)
$ws.Range("H12").Value = 441
$ws.Range("I12").Value = 301.5
$ws.Range("K12").Value = 301.5
$ws.Range("M12").Value = -131.5

# Row 40 (G40=5505)
$ws.Range("H40").Value = 7570.4644
$ws.Range("J40").Value = 8368.608
$ws.Range("L40").Value = 8368.608
$ws.Range("N40").Value = -8718.608

# Row 55 (G55=5517)
$ws.Range("H55").Value = 300

# Row 111 (G111=27768)
$ws.Range("H111").Value = 4810.778
$ws.Range("I111").Value = 4328.143
$ws.Range("K111").Value = 12984.429
$ws.Range("M111").Value = -9917.429

# Row 137 (G137=44013)
$ws.Range("H137").Value = 1999.6666
$ws.Range("I137").Value = 2000
$ws.Range("J137").Value = 1999.5
$ws.Range("K137").Value = 6000
$ws.Range("L137").Value = 5998.5
$ws.Range("M137").Value = -3450
$ws.Range("N137").Value = -11098.5

$ws = $wb.Worksheets.Item("ARM")
# Row 3 (G3=2494)
$ws.Range("H3").Value = 4130
$ws.Range("J3").Value = 6833.3335
$ws.Range("L3").Value = 6833.3335
$ws.Range("N3").Value = -7063.3335

# Row 74 (G74=44000)
$ws.Range("H74").Value = 900
$ws.Range("I74").Value = 900
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 900
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -26
$ws.Range("N74").ClearContents()

# Row 77 (G77=44000)
$ws.Range("H77").Value = 900
$ws.Range("I77").Value = 900
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 4500
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -132
$ws.Range("N77").ClearContents()

# Row 101 (G101=18518)
$ws.Range("H101").Value = 36967
$ws.Range("J101").Value = 36967
$ws.Range("L101").Value = 36967
$ws.Range("N101").Value = -43457

# Row 104 (G104=18672)
$ws.Range("H104").Value = 28870.666
$ws.Range("J104").Value = 28870.666
$ws.Range("L104").Value = 28870.666
$ws.Range("N104").Value = -35858.666

# Row 122 (G122=36168)
$ws.Range("H122").Value = 4948.5
$ws.Range("I122").Value = 4899
$ws.Range("J122").Value = 4998
$ws.Range("K122").Value = 14697
$ws.Range("L122").Value = 14994
$ws.Range("N122").Value = -19894
$ws.Range("M122").Value = -12247

# Row 132 (G132=43997)
$ws.Range("H132").Value = 1406.75
$ws.Range("I132").Value = 1406.75
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 4220.25
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1690.25
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 76 (G76=10630)
$ws.Range("H76").Value = 16999.5
$ws.Range("J76").Value = 16999.5
$ws.Range("L76").Value = 16999.5
$ws.Range("N76").Value = -17629.5

# Row 79 (G79=10630)
$ws.Range("H79").Value = 16999.5
$ws.Range("J79").Value = 16999.5
$ws.Range("L79").Value = 16999.5
$ws.Range("N79").Value = -19183.5

# Row 82 (G82=11877)
$ws.Range("H82").Value = 16052.333
$ws.Range("I82").Value = 16052.333
$ws.Range("K82").Value = 16052.333
$ws.Range("M82").Value = -15669.333

# Row 85 (G85=11877)
$ws.Range("H85").Value = 16052.333
$ws.Range("I85").Value = 16052.333
$ws.Range("K85").Value = 16052.333
$ws.Range("M85").Value = -14726.333

$ws = $wb.Worksheets.Item("CRP")
# Row 4 (G4=3742)
$ws.Range("H4").Value = 1041.6666
$ws.Range("I4").Value = 888.8889
$ws.Range("J4").Value = 1500
$ws.Range("K4").Value = 888.8889
$ws.Range("L4").Value = 1500
$ws.Range("M4").Value = -776.8889
$ws.Range("N4").Value = -1724

# Row 31 (G31=44023)
$ws.Range("H31").Value = 2239.5625
$ws.Range("J31").Value = 3247
$ws.Range("L31").Value = 3247
$ws.Range("N31").Value = -3837

# Row 34 (G34=44023)
$ws.Range("H34").Value = 2239.5625
$ws.Range("J34").Value = 3247
$ws.Range("L34").Value = 3247
$ws.Range("N34").Value = -3651

# Row 94 (G94=32934)
$ws.Range("H94").Value = 1099.25
$ws.Range("J94").Value = 1333.3334
$ws.Range("L94").Value = 1333.3334
$ws.Range("N94").Value = -2235.3334

# Row 106 (G106=18661)
$ws.Range("H106").Value = 30671
$ws.Range("J106").Value = 30671
$ws.Range("L106").Value = 30671
$ws.Range("N106").Value = -33195

# Row 122 (G122=36196)
$ws.Range("H122").Value = 1214.7778
$ws.Range("I122").Value = 1224.4
$ws.Range("K122").Value = 3673.2
$ws.Range("M122").Value = -1223.2

# Row 134 (G134=44020)
$ws.Range("H134").Value = 2126.8572
$ws.Range("I134").Value = 1980
$ws.Range("K134").Value = 5940
$ws.Range("M134").Value = -3405

$ws = $wb.Worksheets.Item("CUL")
# Row 19 (G19=4682)
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()

# Row 39 (G39=4712)
$ws.Range("H39").Value = 2488.3333
$ws.Range("I39").Value = 2178.2
$ws.Range("J39").Value = 2876
$ws.Range("K39").Value = 6534.599999999999
$ws.Range("L39").Value = 8628
$ws.Range("M39").Value = -6240.599999999999
$ws.Range("N39").Value = -9216

# Row 98 (G98=19843)
$ws.Range("H98").Value = 775.44446
$ws.Range("J98").Value = 466
$ws.Range("L98").Value = 1398
$ws.Range("N98").Value = -4394

$ws = $wb.Worksheets.Item("GSM")
# Row 35 (G35=4317)
$ws.Range("H35").Value = 10000
$ws.Range("I35").Value = 10000
$ws.Range("K35").Value = 10000
$ws.Range("M35").Value = -9702

# Row 102 (G102=36169)
$ws.Range("H102").Value = 3763.5715
$ws.Range("I102").Value = 3763.5715
$ws.Range("K102").Value = 3763.5715
$ws.Range("M102").Value = -2141.5715

# Row 105 (G105=18671)
$ws.Range("H105").Value = 20671
$ws.Range("J105").Value = 20671
$ws.Range("L105").Value = 20671
$ws.Range("N105").Value = -27659

# Row 122 (G122=36182)
$ws.Range("H122").Value = 10006.5
$ws.Range("I122").Value = 10006.5
$ws.Range("K122").Value = 30019.5
$ws.Range("M122").Value = -27569.5

# Row 132 (G132=44008)
$ws.Range("H132").Value = 7004
$ws.Range("I132").Value = 5806.8
$ws.Range("J132").Value = 9997
$ws.Range("K132").Value = 17420.4
$ws.Range("L132").Value = 29991
$ws.Range("M132").Value = -14890.4
$ws.Range("N132").Value = -35051

# Row 136 (G136=42218)
$ws.Range("H136").Value = 40000
$ws.Range("J136").Value = 40000
$ws.Range("L136").Value = 120000
$ws.Range("N136").Value = -125100

$ws = $wb.Worksheets.Item("LTW")
# Row 32 (G32=2250)
$ws.Range("H32").Value = 500
$ws.Range("I32").Value = 500
$ws.Range("K32").Value = 500
$ws.Range("M32").Value = -183

# Row 35 (G35=1697)
$ws.Range("H35").Value = 687.3333
$ws.Range("I35").Value = 687.3333
$ws.Range("K35").Value = 687.3333
$ws.Range("M35").Value = -351.3333

# Row 61 (G61=27740)
$ws.Range("H61").Value = 4700
$ws.Range("I61").Value = 5400
$ws.Range("K61").Value = 5400
$ws.Range("M61").Value = -5198

# Row 63 (G63=12006)
$ws.Range("H63").Value = 70000
$ws.Range("J63").Value = 70000
$ws.Range("L63").Value = 70000
$ws.Range("N63").Value = -71498

# Row 66 (G66=12006)
$ws.Range("H66").Value = 70000
$ws.Range("J66").Value = 70000
$ws.Range("L66").Value = 210000
$ws.Range("N66").Value = -217488

# Row 100 (G100=19995)
$ws.Range("H100").Value = 5999.75
$ws.Range("I100").Value = 4666.3335
$ws.Range("J100").Value = 10000
$ws.Range("K100").Value = 4666.3335
$ws.Range("L100").Value = 10000
$ws.Range("M100").Value = -4125.3335
$ws.Range("N100").Value = -11082

# Row 113 (G113=27740)
$ws.Range("H113").Value = 4700
$ws.Range("I113").Value = 5400
$ws.Range("K113").Value = 5400
$ws.Range("M113").Value = -3230

# Row 122 (G122=36247)
$ws.Range("H122").Value = 6150
$ws.Range("I122").Value = 6150
$ws.Range("K122").Value = 18450
$ws.Range("M122").Value = -16000

# Row 128 (G128=34582)
$ws.Range("H128").Value = 99995
$ws.Range("J128").Value = 99995
$ws.Range("L128").Value = 99995
$ws.Range("N128").Value = -109955

# Row 132 (G132=44058)
$ws.Range("H132").Value = 8875
$ws.Range("I132").Value = 7500
$ws.Range("K132").Value = 22500
$ws.Range("M132").Value = -19970

$ws = $wb.Worksheets.Item("WVR")
# Row 101 (G101=18538)
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

# Row 103 (G103=18548)
$ws.Range("H103").Value = 69994.5
$ws.Range("J103").Value = 69994.5
$ws.Range("L103").Value = 69994.5
$ws.Range("N103").Value = -72338.5

# Row 122 (G122=36208)
$ws.Range("H122").Value = 2952.3333
$ws.Range("I122").Value = 3453.2727
$ws.Range("K122").Value = 10359.8181
$ws.Range("M122").Value = -7909.8181

# Row 132 (G132=44029)
$ws.Range("H132").Value = 2398.6365
$ws.Range("I132").Value = 2388.7778
$ws.Range("K132").Value = 7166.3334
$ws.Range("M132").Value = -4636.3334
